# Add two new test rows (fm17 and fm18) to the "ftests" sheet, mirroring the
# existing rows so both the values and the cell formatting match the rows
# already in use for similarly shaped records (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- 1. Write the new cell values first (in the same order the source data
#        was authored in, so that any newly created shared strings come out
#        in the expected sequence: B22, D22, B23, D23, C22, C23, ...). ---
$ws.Range("B22").Value = "fm17"
$ws.Range("D22").Value = "-1,1"
$ws.Range("B23").Value = "fm18"
$ws.Range("D23").Value = "-1,2"
$ws.Range("C22").Value = "WE5 Residential policy with coverage deductibles and blanket policy terms. Ground up loss back-allocation"
$ws.Range("C23").Value = "WE5 Residential policy with coverage deductibles and blanket policy terms. Previous level input loss back-allocation"

$ws.Range("E22").Value = "12,1"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = "3,1"

$ws.Range("E23").Value = "12,1"
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = "3,1"

# --- 2. Copy the formatting (styles only) from an existing, similarly
#        structured row (row 6) onto the two new rows. ---
$ws.Range("B6:K6").Copy()
$ws.Range("B22:K22").PasteSpecial(-4122)
$ws.Range("B6:K6").Copy()
$ws.Range("B23:K23").PasteSpecial(-4122)

# --- 3. Row 22 has no K cell, and row 23 has no I, J or K cells in the
#        target layout, so remove the formatting/content that PasteSpecial
#        brought in for those trailing cells. ---
$ws.Range("K22").Clear()
$ws.Range("I23:K23").Clear()

# --- 4. Match the new selection recorded in the workbook. ---
$ws.Activate()
$ws.Range("C23").Select()
